$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out all existing content on the sheet (old data spanned D2:L10)
$ws.Cells.Clear()

# Write the new values
$ws.Range("I3").Value = "sdlkfjsldfkdslkfj"
$ws.Range("G5").Value = "sdflksjdflskdjflskdjflskdjfldskjf"

# Update the active selection to match the final state
$ws.Range("G5").Select()
